$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose F/G cells previously held the placeholder text "not Found"
# and now get real computed flux values (per commit: "add values for
# unfound results"). Rows 2, 5 and 12 keep the scientific-notation format
# used elsewhere in the sheet (style index 2 / numFmtId 11); the rest use
# the plain General format.

$ws.Range("F2").NumberFormat = "0.00E+00"
$ws.Range("G2").NumberFormat = "0.00E+00"
$ws.Range("F2").Value = [double]"-4.268E-16"
$ws.Range("G2").Value = [double]"3.0820000000000001E-17"

$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

$ws.Range("F5").NumberFormat = "0.00E+00"
$ws.Range("G5").NumberFormat = "0.00E+00"
$ws.Range("F5").Value = [double]"-1.88304E-12"
$ws.Range("G5").Value = [double]"-1.3634E-14"

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

$ws.Range("F8").Value = 0.64849999999999997
$ws.Range("G8").Value = 3.9830999999999999

$ws.Range("F12").NumberFormat = "0.00E+00"
$ws.Range("G12").NumberFormat = "0.00E+00"
$ws.Range("F12").Value = [double]"-9.6799999999999995E-16"
$ws.Range("G12").Value = [double]"9.8422999999999998E-15"

$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

$ws.Range("F22").Value = 2.0880000000000001
$ws.Range("G22").Value = 7.6689999999999996

$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0

$ws.Range("F39").Value = 2.8721999999999999
$ws.Range("G39").Value = 22.842600000000001

$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0

$ws.Range("F41").Value = 2.0886
$ws.Range("G41").Value = 7.6689999999999996

$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0

# Restore the cursor/selection left by the author at F32
$ws.Range("F32").Select()
